$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsETLE = $wb.Worksheets.Item("ETLE")

# Clear existing content first
$wsAbout.Cells.Clear() | Out-Null
$wsETLE.Cells.Clear() | Out-Null

# ---------- write in shared-string creation order ----------

# 0: "Exponent" -> ETLE B1
$wsETLE.Range("B1").Value = "Exponent"

# 1: "ETLE Electricity Technology Logit Exponent" -> About A1
$wsAbout.Range("A1").Value = "ETLE Electricity Technology Logit Exponent"
$wsAbout.Range("A1").Font.Bold = $true

# 2: "Sources:" -> About A3
$wsAbout.Range("A3").Value = "Sources:"
$wsAbout.Range("A3").Font.Bold = $true

# 3: "None needed.  Handled through calibration." -> About B3
$wsAbout.Range("B3").Value = "None needed.  Handled through calibration."

# 4: "Notes" -> About A5
$wsAbout.Range("A5").Value = "Notes"
$wsAbout.Range("A5").Font.Bold = $true

# 5-9: note paragraph -> About A6:A10
$wsAbout.Range("A6").Value = "The logit exponents express how large of a cost difference between technology options"
$wsAbout.Range("A7").Value = "is required to produce a change in technology selection.  This parameter needs to be"
$wsAbout.Range("A8").Value = "obtained via model calibration - e.g. testing a given price intervention with different"
$wsAbout.Range("A9").Value = "logit exponent values until it produces a technology choice shift that matches real-world"
$wsAbout.Range("A10").Value = "data on technology buyers' behavior."

# 10: For more on this... -> About A12
$wsAbout.Range("A12").Value = 'For more on this, see the "Modified Logit" equation description at:'

# 11: url -> About A13
$wsAbout.Range("A13").Value = "https://jgcri.github.io/gcam-doc/choice.html"

# 12: "all electricity sources" -> ETLE A2
$wsETLE.Range("A2").Value = "all electricity sources"
$wsETLE.Range("B2").Value = -3

# 13: "Unit: dimensionless" -> ETLE A1
$wsETLE.Range("A1").Value = "Unit: dimensionless"
$wsETLE.Range("A1").Font.Italic = $true

# ---------- styles ----------
$wsETLE.Range("B1").HorizontalAlignment = -4152  # xlRight

# ---------- column widths ----------
# Column 1 already has width 34.5703125 / bestFit in the source workbook - leave untouched.
# Column 2 is new; closest representable width to the target 12.140625 is used.
$wsETLE.Columns.Item(2).ColumnWidth = 11.33

# ---------- page setup ----------
$wsETLE.PageSetup.Orientation = 1  # xlPortrait

# ---------- selection / active cell ----------
$wsETLE.Range("A1").Select() | Out-Null
$wsAbout.Range("A1").Select() | Out-Null
$wsAbout.Activate() | Out-Null
